# Refresh the "Coeficients" sheet with the latest model run: update the
# label + coefficient for every existing row and append the new rows that
# came out of the expanded model (rows 50-66).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coeficients")

$coeficients = @(
    @{ Row = 2; Label = 'Opportunity Amount < 1500000'; Value = -6.939633873331229 }
    @{ Row = 3; Label = 'Deal Size Category_7'; Value = -5.527185074322019 }
    @{ Row = 4; Label = 'Opportunity Amount < 800000'; Value = -2.236279411676732 }
    @{ Row = 5; Label = 'Deal Size Category_6'; Value = -2.060590304364239 }
    @{ Row = 6; Label = 'Opportunity Amount > 700000'; Value = -1.777044912875483 }
    @{ Row = 7; Label = 'Opportunity Amount < 700000'; Value = -1.681807208628339 }
    @{ Row = 8; Label = 'Opportunity Amount < 500000'; Value = -1.412448799009194 }
    @{ Row = 9; Label = 'POS_Telecoverage'; Value = -1.260737449773282 }
    @{ Row = 10; Label = 'Opportunity Amount > 800000'; Value = -0.9795805875390046 }
    @{ Row = 11; Label = 'Opportunity Amount > 600000'; Value = -0.7712545717361659 }
    @{ Row = 12; Label = 'Deal Size Category_5'; Value = -0.6878411981693952 }
    @{ Row = 13; Label = 'Category_Tires & Wheels'; Value = -0.2581571898719109 }
    @{ Row = 14; Label = 'Region_Northwest'; Value = -0.1621776615663869 }
    @{ Row = 15; Label = 'Deal Size Category_2'; Value = -0.1560359730969999 }
    @{ Row = 16; Label = 'Deal Size Category_4'; Value = -0.1500991996522717 }
    @{ Row = 17; Label = 'Competitor_Unknown'; Value = -0.1391273140944808 }
    @{ Row = 18; Label = 'Client Size By Revenue_5'; Value = -0.1011055924866702 }
    @{ Row = 19; Label = 'Client Size By Employee Count_2'; Value = -0.05541413434382469 }
    @{ Row = 20; Label = 'POS_Other'; Value = -0.05537793870631925 }
    @{ Row = 21; Label = 'Client Size By Employee Count_5'; Value = -0.02743676828376685 }
    @{ Row = 22; Label = 'Opportunity Amount > 1000000'; Value = 0 }
    @{ Row = 23; Label = 'Region_Southwest'; Value = 0.004393777114749361 }
    @{ Row = 24; Label = 'POS_Telesales'; Value = 0.008055282051123522 }
    @{ Row = 25; Label = 'Category_Performance & Non-auto'; Value = 0.008196743063353001 }
    @{ Row = 26; Label = 'Client Size By Employee Count_4'; Value = 0.009548385559475542 }
    @{ Row = 27; Label = 'Client Size By Revenue_2'; Value = 0.01614071704645391 }
    @{ Row = 28; Label = 'Region_Southeast'; Value = 0.01791690768178388 }
    @{ Row = 29; Label = 'Client Size By Employee Count_3'; Value = 0.01950298737272416 }
    @{ Row = 30; Label = 'Client Size By Revenue_3'; Value = 0.02522245168138035 }
    @{ Row = 31; Label = 'Region_Pacific'; Value = 0.02647066573052741 }
    @{ Row = 32; Label = 'Region_Northeast'; Value = 0.05847413352479494 }
    @{ Row = 33; Label = 'Deal Size Category_3'; Value = 0.1033349428439141 }
    @{ Row = 34; Label = 'Client Size By Revenue_4'; Value = 0.1149792450890467 }
    @{ Row = 35; Label = 'Region_Midwest'; Value = 0.1222506332726698 }
    @{ Row = 36; Label = 'POS_Reseller'; Value = 0.3256259525853997 }
    @{ Row = 37; Label = 'Competitor_None'; Value = 0.3493955740015249 }
    @{ Row = 38; Label = 'Category_Car Electronics'; Value = 0.3830181436576884 }
    @{ Row = 39; Label = 'Opportunity Amount < 600000'; Value = 0.61212581016317 }
    @{ Row = 40; Label = 'Opportunity Amount < 250000'; Value = 0.648141505355007 }
    @{ Row = 41; Label = 'Opportunity Amount > 1000'; Value = 0.851826258322658 }
    @{ Row = 42; Label = 'Opportunity Amount > 10000'; Value = 0.9285736608932371 }
    @{ Row = 43; Label = 'Opportunity Amount > 25000'; Value = 0.9417053407178654 }
    @{ Row = 44; Label = 'Opportunity Amount > 100000'; Value = 1.295855238212496 }
    @{ Row = 45; Label = 'Opportunity Amount > 5000'; Value = 1.305965798225456 }
    @{ Row = 46; Label = 'Opportunity Amount < 100000'; Value = 1.335982703524365 }
    @{ Row = 47; Label = 'Opportunity Amount < 25000'; Value = 1.382746960332757 }
    @{ Row = 48; Label = 'Opportunity Amount < 150000'; Value = 1.393470421993614 }
    @{ Row = 49; Label = 'Opportunity Amount < 5000'; Value = 1.46245103607736 }
    @{ Row = 50; Label = 'Opportunity Amount < 50000'; Value = 1.48608190317662 }
    @{ Row = 51; Label = 'Opportunity Amount < 400000'; Value = 1.502331698125935 }
    @{ Row = 52; Label = 'Opportunity Amount < 10000'; Value = 1.538782933429707 }
    @{ Row = 53; Label = 'Opportunity Amount > 50000'; Value = 1.561398754893629 }
    @{ Row = 54; Label = 'Opportunity Amount < 1000000'; Value = 1.583482590393009 }
    @{ Row = 55; Label = 'Revenue From Client Past Two Years_4'; Value = 1.843494631730819 }
    @{ Row = 56; Label = 'Opportunity Amount < 200000'; Value = 1.910307626268267 }
    @{ Row = 57; Label = 'Opportunity Amount > 400000'; Value = 1.932790626035331 }
    @{ Row = 58; Label = 'Opportunity Amount > 300000'; Value = 1.96642652539062 }
    @{ Row = 59; Label = 'Opportunity Amount > 150000'; Value = 2.032158917394971 }
    @{ Row = 60; Label = 'Opportunity Amount > 200000'; Value = 2.181795253150293 }
    @{ Row = 61; Label = 'Revenue From Client Past Two Years_3'; Value = 2.33156511099337 }
    @{ Row = 62; Label = 'Opportunity Amount > 500000'; Value = 2.530815313236761 }
    @{ Row = 63; Label = 'Revenue From Client Past Two Years_2'; Value = 2.674484168553326 }
    @{ Row = 64; Label = 'Opportunity Amount > 250000'; Value = 2.675726861948516 }
    @{ Row = 65; Label = 'Revenue From Client Past Two Years_1'; Value = 2.926876266517279 }
    @{ Row = 66; Label = 'Opportunity Amount < 300000'; Value = 3.095470235370601 }
)

# Rows 50-66 are new: extend column A's formatting (bold font, thin
# border, centered/top alignment) down from the last previously-formatted
# row before writing into them.
$formatSource = $ws.Range("A49")
$formatTarget = $ws.Range("A50:A66")
$formatSource.Copy($formatTarget)

foreach ($item in $coeficients) {
    $ws.Cells.Item($item.Row, 1).Value2 = $item.Label
    $ws.Cells.Item($item.Row, 2).Value2 = $item.Value
}

Write-Output "Coeficients sheet updated through row $($ws.UsedRange.Rows.Count)"
